$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the new bet buttons (bet25, bet50) right after the
# existing "bet10" row (row 3), pushing bet100/spin/exit/... etc down by 2.
$ws.Rows("4:5").Insert()

# --- Update the bet buttons block (rows 2-7) ---
# row 2: bet5
$ws.Range("B2").Value = 48
$ws.Range("C2").Value = 48
$ws.Range("D2").Value = 178
$ws.Range("E2").Value = 307

# row 3: bet10
$ws.Range("B3").Value = 48
$ws.Range("C3").Value = 48
$ws.Range("D3").Value = 237
$ws.Range("E3").Value = 307

# row 4: bet25 (new)
$ws.Range("A4").Value = "bet25"
$ws.Range("B4").Value = 48
$ws.Range("C4").Value = 48
$ws.Range("D4").Value = 298
$ws.Range("E4").Value = 307

# row 5: bet50 (new)
$ws.Range("A5").Value = "bet50"
$ws.Range("B5").Value = 48
$ws.Range("C5").Value = 48
$ws.Range("D5").Value = 357
$ws.Range("E5").Value = 307

# row 6: bet100 (was row 4)
$ws.Range("B6").Value = 48
$ws.Range("C6").Value = 48
$ws.Range("D6").Value = 417
$ws.Range("E6").Value = 307

# row 7: spin (was row 5)
$ws.Range("B7").Value = 131
$ws.Range("C7").Value = 37
$ws.Range("D7").Value = 265
$ws.Range("E7").Value = 401

# rows 8 (exit) and 9 (reset) are unchanged aside from the row shift, which
# the Insert() already handled.

# --- Update the label row block (previously rows 16-19, now 18-21) ---
# row 18: Jackpot txt
$ws.Range("E18").Value = 99

# row 19: Credits
$ws.Range("D19").Value = 223
$ws.Range("E19").Value = 247

# row 20: Bet
$ws.Range("D20").Value = 288
$ws.Range("E20").Value = 247

# row 21: Result
$ws.Range("D21").Value = 355
$ws.Range("E21").Value = 247

# --- View changes ---
[void]$ws.Range("D9").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
